
$wb = $excel.ActiveWorkbook

# ---- Shared literal values for the new "handoff" file row ----
$fname           = '603f375d-b3f6-4fb9-b959-4ee6f197888dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$pathname        = 'e2e\603f375d-b3f6-4fb9-b959-4ee6f197888dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$status          = 'Ready for handoff'
$dateHo          = '2016-09-04 06:30:19'
$zhcnXlf         = '603f375d-b3f6-4fb9-b959-4ee6f197888doooooooooooooooooooooooooooooooooooooooo.274eaf9352a64ff1def5b7ab0aca50f487dded3c.zh-cn.xlf'
$dateHandoffZh   = '2016-09-04 06:30:14'
$dedeXlf         = '603f375d-b3f6-4fb9-b959-4ee6f197888doooooooooooooooooooooooooooooooooooooooo.274eaf9352a64ff1def5b7ab0aca50f487dded3c.de-de.xlf'
$newUrl          = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4652ed5258bb619357eb3da6f8c46bd1637e59f4/e2e/603f375d-b3f6-4fb9-b959-4ee6f197888dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

# =========================================================
# Sheet 1: "Overview"
# =========================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Cells.Item(3,1).Value = $fname
$wsOverview.Cells.Item(3,3).Value = ".md"
$wsOverview.Cells.Item(3,4).Value = ""
$wsOverview.Cells.Item(3,5).Value = $status
$wsOverview.Cells.Item(3,6).Value = $status
$wsOverview.Cells.Item(3,7).Value = $dateHo

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3,2), $newUrl, "", "", $pathname) | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38

# =========================================================
# Sheet 2: "zh-cn"
# =========================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Cells.Item(3,2).Value = ".md"
$wsZhCn.Cells.Item(3,3).Value = $status
$wsZhCn.Cells.Item(3,4).Value = "e2e"
$wsZhCn.Cells.Item(3,5).Value = "ht"
$wsZhCn.Cells.Item(3,6).Value = "False"
$wsZhCn.Cells.Item(3,7).Value = $zhcnXlf
$wsZhCn.Cells.Item(3,8).Value = $dateHandoffZh
$wsZhCn.Cells.Item(3,9).Value = ""
$wsZhCn.Cells.Item(3,10).Value = ""
$wsZhCn.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3,12).Value = ""
$wsZhCn.Cells.Item(3,13).Value = "True"
$wsZhCn.Cells.Item(3,14).Value = ""
$wsZhCn.Cells.Item(3,15).Value = "False"
$wsZhCn.Cells.Item(3,16).Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3,1), $newUrl, "", "", $fname) | Out-Null

$wsZhCn.Columns.Item(3).ColumnWidth = 16.38

# =========================================================
# Sheet 3: "de-de"
# =========================================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Cells.Item(3,2).Value = ".md"
$wsDeDe.Cells.Item(3,3).Value = $status
$wsDeDe.Cells.Item(3,4).Value = "e2e"
$wsDeDe.Cells.Item(3,5).Value = "ht"
$wsDeDe.Cells.Item(3,6).Value = "False"
$wsDeDe.Cells.Item(3,7).Value = $dedeXlf
$wsDeDe.Cells.Item(3,8).Value = $dateHo
$wsDeDe.Cells.Item(3,9).Value = ""
$wsDeDe.Cells.Item(3,10).Value = ""
$wsDeDe.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3,12).Value = ""
$wsDeDe.Cells.Item(3,13).Value = "True"
$wsDeDe.Cells.Item(3,14).Value = ""
$wsDeDe.Cells.Item(3,15).Value = "False"
$wsDeDe.Cells.Item(3,16).Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3,1), $newUrl, "", "", $fname) | Out-Null

$wsDeDe.Columns.Item(3).ColumnWidth = 16.38
